$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3296.6667
$ws.Range("I18").Value = 3296.6667
$ws.Range("K18").Value = 3296.6667
$ws.Range("M18").Value = -3012.6667
$ws.Range("H51").Value = 4021.7827
$ws.Range("I51").Value = 4999.9165
$ws.Range("J51").Value = 2954.7273
$ws.Range("K51").Value = 4999.9165
$ws.Range("L51").Value = 2954.7273
$ws.Range("M51").Value = -4515.9165
$ws.Range("N51").Value = -3922.7273
$ws.Range("H55").Value = 167
$ws.Range("I55").Value = 48.333332
$ws.Range("J55").Value = 226.33333
$ws.Range("K55").Value = 48.333332
$ws.Range("L55").Value = 226.33333
$ws.Range("M55").Value = 165.666668
$ws.Range("N55").Value = -654.3333299999999
$ws.Range("H74").Value = 5135.7144
$ws.Range("I74").Value = 4077.889
$ws.Range("K74").Value = 4077.889
$ws.Range("M74").Value = -3141.889
$ws.Range("H77").Value = 5135.7144
$ws.Range("I77").Value = 4077.889
$ws.Range("K77").Value = 20389.445
$ws.Range("M77").Value = -15709.445
$ws.Range("H96").Value = 2297.6667
$ws.Range("J96").Value = 3207.889
$ws.Range("L96").Value = 9623.667000000001
$ws.Range("N96").Value = -12369.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2595.2727
$ws.Range("I61").Value = 2413.0667
$ws.Range("J61").Value = 2985.7144
$ws.Range("K61").Value = 2413.0667
$ws.Range("L61").Value = 2985.7144
$ws.Range("M61").Value = -2201.0667
$ws.Range("N61").Value = -3409.7144
$ws.Range("H63").Value = 3603.3845
$ws.Range("J63").Value = 5176
$ws.Range("L63").Value = 5176
$ws.Range("N63").Value = -6548
$ws.Range("H66").Value = 3603.3845
$ws.Range("J66").Value = 5176
$ws.Range("L66").Value = 25880
$ws.Range("N66").Value = -32744
$ws.Range("H102").Value = 2348.2856
$ws.Range("I102").Value = 2008.4
$ws.Range("K102").Value = 2008.4
$ws.Range("M102").Value = -386.4000000000001
$ws.Range("H128").Value = 44995
$ws.Range("J128").Value = 44995
$ws.Range("L128").Value = 44995
$ws.Range("N128").Value = -54955
$ws.Range("H132").Value = 2734.0977
$ws.Range("I132").Value = 2340.0967
$ws.Range("J132").Value = 3955.5
$ws.Range("K132").Value = 7020.2901
$ws.Range("L132").Value = 11866.5
$ws.Range("M132").Value = -4490.2901
$ws.Range("N132").Value = -16926.5
$ws.Range("H136").Value = 2595.2727
$ws.Range("I136").Value = 2413.0667
$ws.Range("J136").Value = 2985.7144
$ws.Range("K136").Value = 7239.2001
$ws.Range("L136").Value = 8957.143199999999
$ws.Range("M136").Value = -4689.2001
$ws.Range("N136").Value = -14057.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 31665.666
$ws.Range("I82").Value = 14998.5
$ws.Range("J82").Value = 65000
$ws.Range("K82").Value = 14998.5
$ws.Range("L82").Value = 65000
$ws.Range("M82").Value = -14615.5
$ws.Range("N82").Value = -65766
$ws.Range("H85").Value = 31665.666
$ws.Range("I85").Value = 14998.5
$ws.Range("J85").Value = 65000
$ws.Range("K85").Value = 14998.5
$ws.Range("L85").Value = 65000
$ws.Range("M85").Value = -13672.5
$ws.Range("N85").Value = -67652
$ws.Range("H94").Value = 1049.1818
$ws.Range("I94").Value = 505.4375
$ws.Range("K94").Value = 505.4375
$ws.Range("M94").Value = -54.4375
$ws.Range("H134").Value = 12650789
$ws.Range("I134").Value = 2748731
$ws.Range("J134").Value = 55559708
$ws.Range("K134").Value = 8246193
$ws.Range("L134").Value = 166679124
$ws.Range("M134").Value = -8243658
$ws.Range("N134").Value = -166684194

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2500.1892
$ws.Range("I58").Value = 1923.0385
$ws.Range("K58").Value = 1923.0385
$ws.Range("M58").Value = -1720.0385
$ws.Range("H75").Value = 104638.664
$ws.Range("J75").Value = 104638.664
$ws.Range("L75").Value = 104638.664
$ws.Range("N75").Value = -106634.664
$ws.Range("H78").Value = 104638.664
$ws.Range("J78").Value = 104638.664
$ws.Range("L78").Value = 313915.992
$ws.Range("N78").Value = -323899.992
$ws.Range("H100").Value = 106995
$ws.Range("J100").Value = 106995
$ws.Range("L100").Value = 106995
$ws.Range("N100").Value = -109159
$ws.Range("H136").Value = 2500.1892
$ws.Range("I136").Value = 1923.0385
$ws.Range("K136").Value = 5769.1155
$ws.Range("M136").Value = -3219.1155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2383
$ws.Range("I81").Value = 800
$ws.Range("J81").Value = 4757.5
$ws.Range("K81").Value = 2400
$ws.Range("L81").Value = 14272.5
$ws.Range("M81").Value = -1277
$ws.Range("N81").Value = -16518.5
$ws.Range("H84").Value = 2383
$ws.Range("I84").Value = 800
$ws.Range("J84").Value = 4757.5
$ws.Range("K84").Value = 7200
$ws.Range("L84").Value = 42817.5
$ws.Range("M84").Value = -1584
$ws.Range("N84").Value = -54049.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 134699
$ws.Range("J128").Value = 134699
$ws.Range("L128").Value = 134699
$ws.Range("N128").Value = -144659

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 10618.042
$ws.Range("I82").Value = 1798.3334
$ws.Range("K82").Value = 1798.3334
$ws.Range("M82").Value = -1437.3334
$ws.Range("H85").Value = 10618.042
$ws.Range("I85").Value = 1798.3334
$ws.Range("K85").Value = 1798.3334
$ws.Range("M85").Value = -550.3334
$ws.Range("H93").Value = 1037.0526
$ws.Range("I93").Value = 787.3570999999999
$ws.Range("K93").Value = 787.3570999999999
$ws.Range("M93").Value = 460.6429000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H81").Value = 3782.8125
$ws.Range("I81").Value = 2189.8572
$ws.Range("J81").Value = 5021.778
$ws.Range("K81").Value = 4379.7144
$ws.Range("L81").Value = 10043.556
$ws.Range("M81").Value = -3318.7144
$ws.Range("N81").Value = -12165.556
$ws.Range("H84").Value = 3782.8125
$ws.Range("I84").Value = 2189.8572
$ws.Range("J84").Value = 5021.778
$ws.Range("K84").Value = 21898.572
$ws.Range("L84").Value = 50217.78
$ws.Range("M84").Value = -16594.572
$ws.Range("N84").Value = -60825.78
$ws.Range("H96").Value = 21973.46
$ws.Range("I96").Value = 11832.167
$ws.Range("J96").Value = 30666
$ws.Range("K96").Value = 11832.167
$ws.Range("L96").Value = 30666
$ws.Range("M96").Value = -10459.167
$ws.Range("N96").Value = -33412
$ws.Range("H100").Value = 1849.4117
$ws.Range("J100").Value = 931
$ws.Range("L100").Value = 1862
$ws.Range("N100").Value = -2944
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
